$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.017335002175468
$ws.Range("D2").Value = 1.023230766211161
$ws.Range("E2").Value = 1.044870710006973
$ws.Range("F2").Value = 1.048072478241247
$ws.Range("I2").Value = 1.026990986116645
$ws.Range("J2").Value = 1.022549704457543
$ws.Range("K2").Value = 1.0260628140663
$ws.Range("L2").Value = 1.047640543773365
$ws.Range("M2").Value = 1.050833343735967
$ws.Range("N2").Value = 1.011639532581074
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.018201239408337
$ws.Range("D3").Value = 1.023866998776276
$ws.Range("E3").Value = 1.046109706768693
$ws.Range("F3").Value = 1.049377655204854
$ws.Range("I3").Value = 1.027115038244567
$ws.Range("J3").Value = 1.023052221112168
$ws.Range("K3").Value = 1.026506346915574
$ws.Range("L3").Value = 1.048689751583451
$ws.Range("M3").Value = 1.051949215731737
$ws.Range("N3").Value = 1.011804644630364
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.018761421877567
$ws.Range("D4").Value = 1.024277794435924
$ws.Range("E4").Value = 1.046912097699479
$ws.Range("F4").Value = 1.050222781689925
$ws.Range("I4").Value = 1.027193092046768
$ws.Range("J4").Value = 1.023376440166865
$ws.Range("K4").Value = 1.026791769376637
$ws.Range("L4").Value = 1.049368776464348
$ws.Range("M4").Value = 1.052671302651871
$ws.Range("N4").Value = 1.011911168796187
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.018996842566925
$ws.Range("D5").Value = 1.024450278603126
$ws.Range("E5").Value = 1.047249586262514
$ws.Range("F5").Value = 1.050578215348564
$ws.Range("I5").Value = 1.027225374511432
$ws.Range("J5").Value = 1.02351251539955
$ws.Range("K5").Value = 1.026911383344049
$ws.Range("L5").Value = 1.049654267693184
$ws.Range("M5").Value = 1.052974879336404
$ws.Range("N5").Value = 1.011955876034784
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.019036366020924
$ws.Range("D6").Value = 1.024479226844005
$ws.Range("E6").Value = 1.047306261697851
$ws.Range("F6").Value = 1.050637902617171
$ws.Range("I6").Value = 1.027230763697381
$ws.Range("J6").Value = 1.023535349728454
$ws.Range("K6").Value = 1.026931444879474
$ws.Range("L6").Value = 1.049702204669068
$ws.Range("M6").Value = 1.053025851892378
$ws.Range("N6").Value = 1.011963378140681
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.018764567895781
$ws.Range("D7").Value = 1.024280100021787
$ws.Range("E7").Value = 1.046916606594436
$ws.Range("F7").Value = 1.05022753045037
$ws.Range("I7").Value = 1.027193525495509
$ws.Range("J7").Value = 1.023378259301566
$ws.Range("K7").Value = 1.026793369149547
$ws.Range("L7").Value = 1.049372591096138
$ws.Range("M7").Value = 1.052675359012547
$ws.Range("N7").Value = 1.011911766473577
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.017627819154985
$ws.Range("D8").Value = 1.023445967286118
$ws.Range("E8").Value = 1.045289295330555
$ws.Range("F8").Value = 1.048513447588053
$ws.Range("I8").Value = 1.027033368527652
$ws.Range("J8").Value = 1.022719727087248
$ws.Range("K8").Value = 1.026213033121653
$ws.Range("L8").Value = 1.047995104842159
$ws.Range("M8").Value = 1.051210449586444
$ws.Range("N8").Value = 1.011695397908634
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.015622222177308
$ws.Range("D9").Value = 1.021969369028857
$ws.Range("E9").Value = 1.042426893004294
$ws.Range("F9").Value = 1.045497472347165
$ws.Range("I9").Value = 1.026734211995565
$ws.Range("J9").Value = 1.021552120459853
$ws.Range("K9").Value = 1.025178406270524
$ws.Range("L9").Value = 1.045568646984945
$ws.Range("M9").Value = 1.048629376263786
$ws.Range("N9").Value = 1.011311732390373
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.014283519766652
$ws.Range("D10").Value = 1.020980513067133
$ws.Range("E10").Value = 1.04052198871599
$ws.Range("F10").Value = 1.04348973455169
$ws.Range("I10").Value = 1.026523429902871
$ws.Range("J10").Value = 1.020768924242617
$ws.Range("K10").Value = 1.024480647016815
$ws.Range("L10").Value = 1.043951508179226
$ws.Range("M10").Value = 1.046908784899736
$ws.Range("N10").Value = 1.011054358682973
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.013703466696065
$ws.Range("D11").Value = 1.020551286838808
$ws.Range("E11").Value = 1.039697922615277
$ws.Range("F11").Value = 1.042621032999385
$ws.Range("I11").Value = 1.026429477233563
$ws.Range("J11").Value = 1.020428663239589
$ws.Range("K11").Value = 1.024176622328128
$ws.Range("L11").Value = 1.043251371274518
$ws.Range("M11").Value = 1.046163763589344
$ws.Range("N11").Value = 1.010942537188162
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.013487951725375
$ws.Range("D12").Value = 1.020391697358036
$ws.Range("E12").Value = 1.039391941703506
$ws.Range("F12").Value = 1.042298456250777
$ws.Range("I12").Value = 1.026394176662674
$ws.Range("J12").Value = 1.020302105687245
$ws.Range("K12").Value = 1.024063410927634
$ws.Range("L12").Value = 1.042991322000628
$ws.Range("M12").Value = 1.045887028966154
$ws.Range("N12").Value = 1.010900945277603
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.013534183008054
$ws.Range("D13").Value = 1.020425936871245
$ws.Range("E13").Value = 1.039457570520446
$ws.Range("F13").Value = 1.042367645666682
$ws.Range("I13").Value = 1.02640176695225
$ws.Range("J13").Value = 1.020329260371215
$ws.Range("K13").Value = 1.024087707953564
$ws.Range("L13").Value = 1.043047102900053
$ws.Range("M13").Value = 1.045946389529862
$ws.Range("N13").Value = 1.010909869433799
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.013685653329702
$ws.Range("D14").Value = 1.020538098292522
$ws.Range("E14").Value = 1.0396726278237
$ws.Range("F14").Value = 1.042594366722326
$ws.Range("I14").Value = 1.02642656748168
$ws.Range("J14").Value = 1.020418205406532
$ws.Range("K14").Value = 1.024167270008299
$ws.Range("L14").Value = 1.043229875269116
$ws.Range("M14").Value = 1.04614088863842
$ws.Range("N14").Value = 1.010939100337791
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.013778971600515
$ws.Range("D15").Value = 1.020607184045274
$ws.Range("E15").Value = 1.039805146770198
$ws.Range("F15").Value = 1.042734069938168
$ws.Range("I15").Value = 1.026441794608778
$ws.Range("J15").Value = 1.02047298493984
$ws.Range("K15").Value = 1.0242162533302
$ws.Range("L15").Value = 1.043342489001847
$ws.Range("M15").Value = 1.046260725843893
$ws.Range("N15").Value = 1.010957102991337
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.014322007893144
$ws.Range("D16").Value = 1.021008977487912
$ws.Range("E16").Value = 1.040576695342614
$ws.Range("F16").Value = 1.043547401210495
$ws.Range("I16").Value = 1.026529608769462
$ws.Range("J16").Value = 1.020791482447698
$ws.Range("K16").Value = 1.024500784384741
$ws.Range("L16").Value = 1.043997975829923
$ws.Range("M16").Value = 1.046958229515023
$ws.Range("N16").Value = 1.011061771982632
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.014662537251999
$ws.Range("D17").Value = 1.021260732999989
$ws.Range("E17").Value = 1.04106087229559
$ws.Range("F17").Value = 1.044057758442936
$ws.Range("I17").Value = 1.026583974531194
$ws.Range("J17").Value = 1.020990964901791
$ws.Range("K17").Value = 1.024678757747788
$ws.Range("L17").Value = 1.044409169735575
$ws.Range("M17").Value = 1.047395755749434
$ws.Range("N17").Value = 1.011127327295363
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.014861124975497
$ws.Range("D18").Value = 1.021407476812859
$ws.Range("E18").Value = 1.04134335893502
$ws.Range("F18").Value = 1.044355505191653
$ws.Range("I18").Value = 1.026615426280833
$ws.Range("J18").Value = 1.021107210394284
$ws.Range("K18").Value = 1.024782384150236
$ws.Range("L18").Value = 1.044649021438558
$ws.Range("M18").Value = 1.0476509580862
$ws.Range("N18").Value = 1.011165528213156
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.014928831914613
$ws.Range("D19").Value = 1.021457495504778
$ws.Range("E19").Value = 1.041439692254396
$ws.Range("F19").Value = 1.044457040050228
$ws.Range("I19").Value = 1.026626106583734
$ws.Range("J19").Value = 1.021146828554858
$ws.Range("K19").Value = 1.024817687101217
$ws.Range("L19").Value = 1.044730806286534
$ws.Range("M19").Value = 1.047737975648664
$ws.Range("N19").Value = 1.011178547560658
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.014626005527375
$ws.Range("D20").Value = 1.021233732445336
$ws.Range("E20").Value = 1.041008917006163
$ws.Range("F20").Value = 1.044002995333756
$ws.Range("I20").Value = 1.026578168371686
$ws.Range("J20").Value = 1.020969573618062
$ws.Range("K20").Value = 1.024659681755186
$ws.Range("L20").Value = 1.044365051574409
$ws.Range("M20").Value = 1.047348813264433
$ws.Range("N20").Value = 1.011120297592858
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.013641050676164
$ws.Range("D21").Value = 1.020505073838866
$ws.Range("E21").Value = 1.039609295677165
$ws.Range("F21").Value = 1.042527600340107
$ws.Range("I21").Value = 1.026419275445858
$ws.Range("J21").Value = 1.020392017983957
$ws.Range("K21").Value = 1.024143848773796
$ws.Range("L21").Value = 1.04317605302216
$ws.Range("M21").Value = 1.046083613525205
$ws.Range("N21").Value = 1.010930494121171
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.013021438453592
$ws.Range("D22").Value = 1.020046036754221
$ws.Range("E22").Value = 1.038729956487271
$ws.Range("F22").Value = 1.041600525245098
$ws.Range("I22").Value = 1.026317045513504
$ws.Range("J22").Value = 1.020027905640231
$ws.Range("K22").Value = 1.023817886646528
$ws.Range("L22").Value = 1.04242855543974
$ws.Range("M22").Value = 1.045288128408859
$ws.Range("N22").Value = 1.010810830768065
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.013349937898562
$ws.Range("D23").Value = 1.020289465937784
$ws.Range("E23").Value = 1.039196049018627
$ws.Range("F23").Value = 1.042091932405339
$ws.Range("I23").Value = 1.026371459947983
$ws.Range("J23").Value = 1.020221021200999
$ws.Range("K23").Value = 1.023990840250929
$ws.Range("L23").Value = 1.042824811635418
$ws.Range("M23").Value = 1.045709830968701
$ws.Range("N23").Value = 1.010874297443649
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.014642512756412
$ws.Range("D24").Value = 1.021245933145648
$ws.Range("E24").Value = 1.041032393134648
$ws.Range("F24").Value = 1.044027740227249
$ws.Range("I24").Value = 1.026580792725164
$ws.Range("J24").Value = 1.020979239755943
$ws.Range("K24").Value = 1.024668301938886
$ws.Range("L24").Value = 1.044384986643367
$ws.Range("M24").Value = 1.047370024550225
$ws.Range("N24").Value = 1.011123474125256
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.016141008876762
$ws.Range("D25").Value = 1.022351896542661
$ws.Range("E25").Value = 1.043166293554335
$ws.Range("F25").Value = 1.046276654696619
$ws.Range("I25").Value = 1.026813553893807
$ws.Range("J25").Value = 1.021854822346624
$ws.Range("K25").Value = 1.025447297955395
$ws.Range("L25").Value = 1.046195851211313
$ws.Range("M25").Value = 1.049296618785497
$ws.Range("N25").Value = 1.011411201603874
